$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 922.7143
$ws.Range("J17").Value = 922.7143
$ws.Range("L17").Value = 2768.1429
$ws.Range("N17").Value = -3104.1429
$ws.Range("H62").Value = 2818.0908
$ws.Range("I62").Value = 2600.6
$ws.Range("J62").Value = 2999.3333
$ws.Range("K62").Value = 2600.6
$ws.Range("L62").Value = 2999.3333
$ws.Range("M62").Value = -1976.6
$ws.Range("N62").Value = -4247.3333
$ws.Range("H65").Value = 2818.0908
$ws.Range("I65").Value = 2600.6
$ws.Range("J65").Value = 2999.3333
$ws.Range("K65").Value = 13003
$ws.Range("L65").Value = 14996.6665
$ws.Range("M65").Value = -9883
$ws.Range("N65").Value = -21236.6665
$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").Value = $null
$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").Value = $null
$ws.Range("H112").Value = 4212.533
$ws.Range("J112").Value = 4212.533
$ws.Range("L112").Value = 12637.599
$ws.Range("N112").Value = -14853.599
$ws.Range("H135").Value = 533.9524
$ws.Range("I135").Value = 517.05554
$ws.Range("K135").Value = 4653.49986
$ws.Range("M135").Value = -2118.49986
$ws.Range("H137").Value = 35928.277
$ws.Range("I137").Value = 1328.0454
$ws.Range("K137").Value = 3984.1362
$ws.Range("M137").Value = -1434.1362
$ws.Range("H138").Value = 2794.0925
$ws.Range("I138").Value = 2428.6553
$ws.Range("J138").Value = 3218
$ws.Range("K138").Value = 7285.965899999999
$ws.Range("L138").Value = 9654
$ws.Range("M138").Value = -2145.965899999999
$ws.Range("N138").Value = -19934

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 343.625
$ws.Range("I5").Value = 222
$ws.Range("J5").Value = 546.3333
$ws.Range("K5").Value = 222
$ws.Range("L5").Value = 546.3333
$ws.Range("M5").Value = -110
$ws.Range("N5").Value = -770.3333
$ws.Range("H32").Value = 3739.3872
$ws.Range("I32").Value = 3071.625
$ws.Range("J32").Value = 6028.857
$ws.Range("K32").Value = 3071.625
$ws.Range("L32").Value = 6028.857
$ws.Range("M32").Value = -2784.625
$ws.Range("N32").Value = -6602.857
$ws.Range("H61").Value = 8510.333000000001
$ws.Range("I61").Value = 3799
$ws.Range("K61").Value = 3799
$ws.Range("M61").Value = -3587
$ws.Range("H74").Value = 1140.9048
$ws.Range("I74").Value = 923.7353000000001
$ws.Range("J74").Value = 2063.875
$ws.Range("K74").Value = 923.7353000000001
$ws.Range("L74").Value = 2063.875
$ws.Range("M74").Value = -49.73530000000005
$ws.Range("N74").Value = -3811.875
$ws.Range("H77").Value = 1140.9048
$ws.Range("I77").Value = 923.7353000000001
$ws.Range("J77").Value = 2063.875
$ws.Range("K77").Value = 4618.6765
$ws.Range("L77").Value = 10319.375
$ws.Range("M77").Value = -250.6765000000005
$ws.Range("N77").Value = -19055.375
$ws.Range("H97").Value = 876.4545000000001
$ws.Range("I97").Value = 757
$ws.Range("J97").Value = 1019.8
$ws.Range("K97").Value = 757
$ws.Range("L97").Value = 1019.8
$ws.Range("M97").Value = -261
$ws.Range("N97").Value = -2011.8
$ws.Range("H110").Value = 2730.4285
$ws.Range("I110").Value = 1620
$ws.Range("K110").Value = 1620
$ws.Range("M110").Value = 425
$ws.Range("H122").Value = 1541.375
$ws.Range("I122").Value = 1541.375
$ws.Range("K122").Value = 4624.125
$ws.Range("M122").Value = -2174.125
$ws.Range("H132").Value = 2201.12
$ws.Range("I132").Value = 1845.8334
$ws.Range("K132").Value = 5537.5002
$ws.Range("M132").Value = -3007.5002
$ws.Range("H136").Value = 8510.333000000001
$ws.Range("I136").Value = 3799
$ws.Range("K136").Value = 11397
$ws.Range("M136").Value = -8847

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 343.625
$ws.Range("I4").Value = 222
$ws.Range("J4").Value = 546.3333
$ws.Range("K4").Value = 222
$ws.Range("L4").Value = 546.3333
$ws.Range("M4").Value = -107
$ws.Range("N4").Value = -776.3333
$ws.Range("H105").Value = 2287.0476
$ws.Range("I105").Value = 2251.611
$ws.Range("K105").Value = 2251.611
$ws.Range("M105").Value = -504.6109999999999
$ws.Range("H107").Value = 1467.44
$ws.Range("J107").Value = 1742.7142
$ws.Range("L107").Value = 1742.7142
$ws.Range("N107").Value = -5582.7142
$ws.Range("H134").Value = 9092.947
$ws.Range("I134").Value = 11634
$ws.Range("J134").Value = 1978
$ws.Range("K134").Value = 34902
$ws.Range("L134").Value = 5934
$ws.Range("M134").Value = -32367
$ws.Range("N134").Value = -11004

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1698.8928
$ws.Range("I31").Value = 1059.4166
$ws.Range("J31").Value = 2178.5
$ws.Range("K31").Value = 1059.4166
$ws.Range("L31").Value = 2178.5
$ws.Range("M31").Value = -764.4166
$ws.Range("N31").Value = -2768.5
$ws.Range("H34").Value = 1698.8928
$ws.Range("I34").Value = 1059.4166
$ws.Range("J34").Value = 2178.5
$ws.Range("K34").Value = 1059.4166
$ws.Range("L34").Value = 2178.5
$ws.Range("M34").Value = -857.4166
$ws.Range("N34").Value = -2582.5
$ws.Range("H58").Value = 2290756.8
$ws.Range("I58").Value = 3954387
$ws.Range("J58").Value = 3265.125
$ws.Range("K58").Value = 3954387
$ws.Range("L58").Value = 3265.125
$ws.Range("M58").Value = -3954184
$ws.Range("N58").Value = -3671.125
$ws.Range("H86").Value = 2300
$ws.Range("I86").Value = 1750.5
$ws.Range("K86").Value = 1750.5
$ws.Range("M86").Value = -627.5
$ws.Range("H89").Value = 2300
$ws.Range("I89").Value = 1750.5
$ws.Range("K89").Value = 8752.5
$ws.Range("M89").Value = -3136.5
$ws.Range("H107").Value = 375.2069
$ws.Range("I107").Value = 399.5
$ws.Range("J107").Value = 335.45456
$ws.Range("K107").Value = 399.5
$ws.Range("L107").Value = 335.45456
$ws.Range("M107").Value = 1520.5
$ws.Range("N107").Value = -4175.45456
$ws.Range("H132").Value = 2326.6956
$ws.Range("I132").Value = 1073.5333
$ws.Range("K132").Value = 3220.5999
$ws.Range("M132").Value = -690.5999000000002
$ws.Range("H136").Value = 2290756.8
$ws.Range("I136").Value = 3954387
$ws.Range("J136").Value = 3265.125
$ws.Range("K136").Value = 11863161
$ws.Range("L136").Value = 9795.375
$ws.Range("M136").Value = -11860611
$ws.Range("N136").Value = -14895.375

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 692.5
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = $null
$ws.Range("H87").Value = 12161
$ws.Range("I87").Value = 6158.1665
$ws.Range("K87").Value = 18474.4995
$ws.Range("M87").Value = -17226.4995
$ws.Range("H90").Value = 12161
$ws.Range("I90").Value = 6158.1665
$ws.Range("K90").Value = 55423.4985
$ws.Range("M90").Value = -49183.4985
$ws.Range("H131").Value = 785.3099999999999
$ws.Range("I131").Value = 525
$ws.Range("J131").Value = 796.15625
$ws.Range("K131").Value = 1575
$ws.Range("L131").Value = 2388.46875
$ws.Range("M131").Value = 3465
$ws.Range("N131").Value = -12468.46875

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3872.25
$ws.Range("I80").Value = 3750
$ws.Range("J80").Value = 3994.5
$ws.Range("K80").Value = 3750
$ws.Range("L80").Value = 3994.5
$ws.Range("M80").Value = -2752
$ws.Range("N80").Value = -5990.5
$ws.Range("H83").Value = 3872.25
$ws.Range("I83").Value = 3750
$ws.Range("J83").Value = 3994.5
$ws.Range("K83").Value = 18750
$ws.Range("L83").Value = 19972.5
$ws.Range("M83").Value = -13758
$ws.Range("N83").Value = -29956.5
$ws.Range("H132").Value = 2566775.2
$ws.Range("I132").Value = 5496350
$ws.Range("K132").Value = 16489050
$ws.Range("M132").Value = -16486520

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 20000
$ws.Range("J38").Value = 20000
$ws.Range("L38").Value = 20000
$ws.Range("N38").Value = -20820
$ws.Range("H46").Value = 2454.5833
$ws.Range("I46").Value = 1419.5
$ws.Range("J46").Value = 2972.125
$ws.Range("K46").Value = 1419.5
$ws.Range("L46").Value = 2972.125
$ws.Range("M46").Value = -1231.5
$ws.Range("N46").Value = -3348.125
$ws.Range("H100").Value = 1096.6666
$ws.Range("I100").Value = 1096.6666
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1096.6666
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = $null
$ws.Range("N100").Value = -555.6666
$ws.Range("H122").Value = 13333.333
$ws.Range("I122").Value = 10000
$ws.Range("K122").Value = 30000
$ws.Range("M122").Value = -27550
$ws.Range("H132").Value = 1623.4103
$ws.Range("I132").Value = 1297.8948
$ws.Range("J132").Value = 1932.65
$ws.Range("K132").Value = 3893.6844
$ws.Range("L132").Value = 5797.950000000001
$ws.Range("M132").Value = -1363.6844
$ws.Range("N132").Value = -10857.95
$ws.Range("H136").Value = 2577.12
$ws.Range("I136").Value = 1601.6471
$ws.Range("J136").Value = 4650
$ws.Range("K136").Value = 4804.9413
$ws.Range("L136").Value = 13950
$ws.Range("M136").Value = -2254.9413
$ws.Range("N136").Value = -19050

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 39996
$ws.Range("J97").Value = 39996
$ws.Range("L97").Value = 39996
$ws.Range("N97").Value = -41978
$ws.Range("H132").Value = 1176.5
$ws.Range("I132").Value = 965.3182
$ws.Range("J132").Value = 3499.5
$ws.Range("K132").Value = 2895.9546
$ws.Range("L132").Value = 10498.5
$ws.Range("M132").Value = -365.9546
$ws.Range("N132").Value = -15558.5
$ws.Range("H136").Value = 22224868
$ws.Range("J136").Value = 2715
$ws.Range("L136").Value = 8145
$ws.Range("N136").Value = -13245
